# Auto-generated script applying cell value updates described by the commit diff.
# Each sheet's cells are updated to reflect refreshed market-price derived figures.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 1803.1538
$ws.Range("I12").Value = 115.75
$ws.Range("J12").Value = 4503
$ws.Range("K12").Value = 115.75
$ws.Range("L12").Value = 4503
$ws.Range("M12").Value = 54.25
$ws.Range("N12").Value = -4843
$ws.Range("H17").Value = 68045.734
$ws.Range("J17").Value = 70382.484
$ws.Range("L17").Value = 211147.452
$ws.Range("N17").Value = -211483.452
$ws.Range("H64").Value = 7587.75
$ws.Range("J64").Value = 12333.333
$ws.Range("L64").Value = 12333.333
$ws.Range("N64").Value = -12829.333
$ws.Range("H67").Value = 7587.75
$ws.Range("J67").Value = 12333.333
$ws.Range("L67").Value = 12333.333
$ws.Range("N67").Value = -14049.333
$ws.Range("H76").Value = 43531810
$ws.Range("I76").Value = 92482.75
$ws.Range("K76").Value = 92482.75
$ws.Range("M76").Value = -92167.75
$ws.Range("H79").Value = 43531810
$ws.Range("I79").Value = 92482.75
$ws.Range("K79").Value = 92482.75
$ws.Range("M79").Value = -91390.75
$ws.Range("H86").Value = 5268257.5
$ws.Range("J86").Value = 5853064
$ws.Range("L86").Value = 5853064
$ws.Range("N86").Value = -5855310
$ws.Range("H89").Value = 5268257.5
$ws.Range("J89").Value = 5853064
$ws.Range("L89").Value = 29265320
$ws.Range("N89").Value = -29276552
$ws.Range("H127").Value = 19446.834
$ws.Range("J127").Value = 1512.5
$ws.Range("L127").Value = 4537.5
$ws.Range("N127").Value = -14457.5
$ws.Range("H135").Value = 2017.4166
$ws.Range("I135").Value = 1343.1666
$ws.Range("J135").Value = 2691.6667
$ws.Range("K135").Value = 12088.4994
$ws.Range("L135").Value = 24225.0003
$ws.Range("M135").Value = -9553.499400000001
$ws.Range("N135").Value = -29295.0003
$ws.Range("H137").Value = 3385.4583
$ws.Range("I137").Value = 1592.25
$ws.Range("K137").Value = 4776.75
$ws.Range("M137").Value = -2226.75
$ws.Range("H138").Value = 7214.2563
$ws.Range("J138").Value = 7631.353
$ws.Range("L138").Value = 22894.059
$ws.Range("N138").Value = -33174.059
$ws.Range("H141").Value = 3286.1177
$ws.Range("I141").Value = 3497.4285
$ws.Range("K141").Value = 10492.2855
$ws.Range("M141").Value = -5312.2855
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2181.7192
$ws.Range("I32").Value = 2181.7192
$ws.Range("K32").Value = 2181.7192
$ws.Range("M32").Value = -1894.7192
$ws.Range("H37").Value = 50523.15
$ws.Range("J37").Value = 49468.277
$ws.Range("L37").Value = 49468.277
$ws.Range("N37").Value = -50014.277
$ws.Range("H44").Value = 29933.334
$ws.Range("J44").Value = 29933.334
$ws.Range("L44").Value = 29933.334
$ws.Range("N44").Value = -30909.334
$ws.Range("H45").Value = 5602.7085
$ws.Range("I45").Value = 2172.875
$ws.Range("J45").Value = 12462.375
$ws.Range("K45").Value = 2172.875
$ws.Range("L45").Value = 12462.375
$ws.Range("M45").Value = -1795.875
$ws.Range("N45").Value = -13216.375
$ws.Range("I61").Value = 4514.091
$ws.Range("J61").Value = 8025
$ws.Range("K61").Value = 4514.091
$ws.Range("L61").Value = 8025
$ws.Range("M61").Value = -4302.091
$ws.Range("N61").Value = -8449
$ws.Range("H74").Value = 70497.5
$ws.Range("I74").Value = 260000
$ws.Range("J74").Value = 7330
$ws.Range("K74").Value = 260000
$ws.Range("L74").Value = 7330
$ws.Range("M74").Value = -259126
$ws.Range("N74").Value = -9078
$ws.Range("H77").Value = 70497.5
$ws.Range("I77").Value = 260000
$ws.Range("J77").Value = 7330
$ws.Range("K77").Value = 1300000
$ws.Range("L77").Value = 36650
$ws.Range("M77").Value = -1295632
$ws.Range("N77").Value = -45386
$ws.Range("H110").Value = 174652.23
$ws.Range("I110").Value = 194161.66
$ws.Range("K110").Value = 194161.66
$ws.Range("M110").Value = -192116.66
$ws.Range("H122").Value = 4685.1934
$ws.Range("I122").Value = 4159.143
$ws.Range("K122").Value = 12477.429
$ws.Range("M122").Value = -10027.429
$ws.Range("I136").Value = 4514.091
$ws.Range("J136").Value = 8025
$ws.Range("K136").Value = 13542.273
$ws.Range("L136").Value = 24075
$ws.Range("M136").Value = -10992.273
$ws.Range("N136").Value = -29175
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1913.3
$ws.Range("I99").Value = 2152.5
$ws.Range("K99").Value = 2152.5
$ws.Range("M99").Value = -654.5
$ws.Range("H134").Value = 37086.312
$ws.Range("I134").Value = 5574.476
$ws.Range("J134").Value = 97245.27
$ws.Range("K134").Value = 16723.428
$ws.Range("L134").Value = 291735.81
$ws.Range("M134").Value = -14188.428
$ws.Range("N134").Value = -296805.81
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H51").Value = 27771.053
$ws.Range("J51").Value = 77662.5
$ws.Range("L51").Value = 77662.5
$ws.Range("N51").Value = -79134.5
$ws.Range("H60").Value = 85160
$ws.Range("J60").Value = 98950
$ws.Range("L60").Value = 98950
$ws.Range("N60").Value = -99972
$ws.Range("H61").Value = 27771.053
$ws.Range("J61").Value = 77662.5
$ws.Range("L61").Value = 77662.5
$ws.Range("N61").Value = -78358.5
$ws.Range("H122").Value = 3250
$ws.Range("I122").Value = 1432.5555
$ws.Range("K122").Value = 4297.666499999999
$ws.Range("M122").Value = -1847.666499999999
$ws.Range("H132").Value = 4403.533
$ws.Range("I132").Value = 4157.923
$ws.Range("K132").Value = 12473.769
$ws.Range("M132").Value = -9943.769
$ws.Range("H134").Value = 504180.44
$ws.Range("I134").Value = 4067.4167
$ws.Range("K134").Value = 12202.2501
$ws.Range("M134").Value = -9667.250100000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 29
$ws.Range("I38").Value = 30.4
$ws.Range("J38").Value = 27.833334
$ws.Range("K38").Value = 91.19999999999999
$ws.Range("L38").Value = 83.50000199999999
$ws.Range("M38").Value = 255.8
$ws.Range("N38").Value = -777.500002
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 695979.5600000001
$ws.Range("I122").Value = 923889.4399999999
$ws.Range("K122").Value = 2771668.32
$ws.Range("M122").Value = -2769218.32
$ws.Range("H126").Value = 5833.222
$ws.Range("I126").Value = 4999
$ws.Range("J126").Value = 5937.5
$ws.Range("K126").Value = 14997
$ws.Range("L126").Value = 17812.5
$ws.Range("M126").Value = -12527  # new cell added
$ws.Range("N126").Value = -22752.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 533107.6
$ws.Range("I7").Value = 839220.4399999999
$ws.Range("K7").Value = 839220.4399999999
$ws.Range("M7").Value = -839108.4399999999
$ws.Range("H46").Value = 5163.2354
$ws.Range("I46").Value = 4108.3335
$ws.Range("K46").Value = 4108.3335
$ws.Range("M46").Value = -3920.3335
$ws.Range("H82").Value = 4863.25
$ws.Range("I82").Value = 2700
$ws.Range("K82").Value = 2700
$ws.Range("M82").Value = -2339
$ws.Range("H85").Value = 4863.25
$ws.Range("I85").Value = 2700
$ws.Range("K85").Value = 2700
$ws.Range("M85").Value = -1452
$ws.Range("H126").Value = 533107.6
$ws.Range("I126").Value = 839220.4399999999
$ws.Range("K126").Value = 2517661.32
$ws.Range("M126").Value = -2515191.32
$ws.Range("H132").Value = 6366.4443
$ws.Range("I132").Value = 4459.6
$ws.Range("J132").Value = 8750
$ws.Range("K132").Value = 13378.8
$ws.Range("L132").Value = 26250
$ws.Range("M132").Value = -10848.8
$ws.Range("N132").Value = -31310
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H4").Value = 5000
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("H122").Value = 43483344
$ws.Range("I122").Value = 76926610
$ws.Range("J122").Value = 7100
$ws.Range("K122").Value = 230779830
$ws.Range("L122").Value = 21300
$ws.Range("M122").Value = -230777380
$ws.Range("N122").Value = -26200
$ws.Range("H126").Value = 3182.5715
$ws.Range("I126").Value = 1450.6154
$ws.Range("K126").Value = 4351.8462
$ws.Range("M126").Value = -1881.8462
$ws.Range("H132").Value = 52228.145
$ws.Range("I132").Value = 3214.6924
$ws.Range("K132").Value = 9644.0772
$ws.Range("M132").Value = -7114.0772

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("N4").ClearContents()
